# Added EquipmentManufacturer and EquipmentManufacturerModel to Equipment.dtsx
# Populates a sample data row on each sheet of the Equipment Manufacturer
# template, formats the used range as Text, and makes "Manufacturer" the
# active sheet/tab (it had been "Model").

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Manufacturer")
$ws2 = $wb.Worksheets.Item("Model")

# Format the header + new data row as Text (@) on both sheets so the IDs /
# codes below (e.g. "534:02:00") are stored as literal text, not parsed as
# numbers/dates.
$ws1.Range("A1:E2").NumberFormat = "@"
$ws2.Range("A1:E2").NumberFormat = "@"

# Enter the new sample row's values. Order matters here: it controls the
# order new entries land in the shared-string table.
$ws1.Range("E2").Value = "[GROUPROW]"
$ws2.Range("B2").Value = "[KEY]"

$ws1.Range("A2").Value = "2134"
$ws1.Range("B2").Value = "532:02:00"
$ws1.Range("C2").Value = "534:02:00"
$ws1.Range("D2").Value = "534:14:00"

$ws2.Range("C2").Value = "534:06:00"
$ws2.Range("D2").Value = "534:07:00"
$ws2.Range("E2").Value = "534:16:00"
$ws2.Range("A2").Value = "2134"

# Reset each sheet's lingering cell selection back to A1 (previously D2 on
# Manufacturer, A2 on Model), then make "Manufacturer" the active tab again
# (it had been "Model").
[void]$ws1.Range("A1").Select()
[void]$ws2.Range("A1").Select()
[void]$ws1.Activate()
[void]$ws1.Range("A1").Select()
